# Refresh the cryptocurrency Price (column D) / Volume(1h) (column E)
# snapshot values pulled in by the scheduled GitHub Actions job.
#
# All of these values are stored as literal text (matching the source
# data feed / the sheet's existing inline-string cells), so any new
# value that happens to look like a plain number is entered with a
# leading apostrophe (Excel's standard 'force text' input) and the
# cell style is then reset to Normal so no stray formatting is left
# behind from the apostrophe/text entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '76.333.21'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '3.031.77'
$ws.Range("E3").Value = '  +3.78%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '''197.79'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").Value = '''617.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.66%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -1.01%  '
$ws.Range("D9").Value = '''0.204'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.19%  '
$ws.Range("D10").Value = '3.029.87'
$ws.Range("E10").Value = '  +3.52%  '
$ws.Range("D11").Value = '''0.434'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.28%  '
$ws.Range("E12").Value = '  -0.73%  '
$ws.Range("D13").Value = '''5.27'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.01%  '
$ws.Range("D14").Value = '3.592.97'
$ws.Range("E14").Value = '  +3.75%  '
$ws.Range("D15").Value = '''28.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").Value = '76.206.86'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").Value = '3.029.07'
$ws.Range("E18").Value = '  +3.21%  '
$ws.Range("D19").Value = '''13.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.89%  '
$ws.Range("D20").Value = '''8.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.95%  '
$ws.Range("D21").Value = '''382.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.52%  '
$ws.Range("D22").Value = '''2.34'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.73%  '
$ws.Range("E23").Value = '  +0.28%  '
$ws.Range("D24").Value = '3.188.80'
$ws.Range("E24").Value = '  +3.18%  '
$ws.Range("D25").Value = '''72.72'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.09%  '
$ws.Range("D26").Value = '''0.997'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.29%  '
$ws.Range("D27").Value = '''4.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.30%  '
$ws.Range("D28").Value = '''9.68'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.13%  '
$ws.Range("E29").Value = '  -0.44%  '
$ws.Range("D30").Value = '''0.995'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("D31").Value = '''8.25'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.72%  '
$ws.Range("E32").Value = '  +0.68%  '
$ws.Range("D33").Value = '''489.54'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.06%  '
$ws.Range("D34").Value = '''1.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.74%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = '''20.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.34%  '
$ws.Range("D37").Value = '''162.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.84%  '
$ws.Range("E38").Value = '  +2.04%  '
$ws.Range("E39").Value = '  +4.21%  '
$ws.Range("D40").Value = '''0.380'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.74%  '
$ws.Range("D41").Value = '''190.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.39%  '
$ws.Range("E42").Value = '  -4.96%  '
$ws.Range("D44").Value = '''0.788'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +19.95%  '
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("D46").Value = '''41.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.47%  '
$ws.Range("E47").Value = '  +4.82%  '
$ws.Range("E48").Value = '  -1.57%  '
$ws.Range("D49").Value = '''2.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.91%  '
$ws.Range("D50").Value = '''0.598'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.61%  '
$ws.Range("D51").Value = '''3.87'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.75%  '
